# Update "想去人数" (want-to-go count) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 2160
$ws1.Range("F12").Value = 1368
$ws1.Range("F21").Value = 59
$ws1.Range("F22").Value = 20
$ws1.Range("F23").Value = 1177

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 2160
$ws4.Range("F13").Value = 1368
$ws4.Range("F22").Value = 59
$ws4.Range("F23").Value = 20
$ws4.Range("F24").Value = 1177
